$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14) -- shifts the existing
# "Late" / "heading" / "Outstanding" columns one to the right.
$ws.Columns.Item(14).Insert()

# Match the new column's stored width (11) -- ColumnWidth needs the
# ~5/6 padding offset subtracted to land on an exact stored width of 11.
$ws.Columns.Item(14).ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active sheet/tab, with the cursor on K17,
# which also clears tabSelected on the previously active "Transactions" sheet.
$ws.Activate()
$ws.Range("K17").Select()
